$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

# Update the scraped_at timestamps in column K (rows 2-48) to reflect the new run time
$ws.Cells.Item(2, 11).Value = "2025-11-04T07:01:19.713733+00:00"
$ws.Cells.Item(3, 11).Value = "2025-11-04T07:01:22.038322+00:00"
$ws.Cells.Item(4, 11).Value = "2025-11-04T07:01:22.038357+00:00"
$ws.Cells.Item(5, 11).Value = "2025-11-04T07:01:22.038377+00:00"
$ws.Cells.Item(6, 11).Value = "2025-11-04T07:01:22.038394+00:00"
$ws.Cells.Item(7, 11).Value = "2025-11-04T07:01:24.391853+00:00"
$ws.Cells.Item(8, 11).Value = "2025-11-04T07:01:24.391883+00:00"
$ws.Cells.Item(9, 11).Value = "2025-11-04T07:01:24.391903+00:00"
$ws.Cells.Item(10, 11).Value = "2025-11-04T07:01:27.199711+00:00"
$ws.Cells.Item(11, 11).Value = "2025-11-04T07:01:29.652989+00:00"
$ws.Cells.Item(12, 11).Value = "2025-11-04T07:01:29.653021+00:00"
$ws.Cells.Item(13, 11).Value = "2025-11-04T07:01:29.653040+00:00"
$ws.Cells.Item(14, 11).Value = "2025-11-04T07:01:32.034472+00:00"
$ws.Cells.Item(15, 11).Value = "2025-11-04T07:01:32.034501+00:00"
$ws.Cells.Item(16, 11).Value = "2025-11-04T07:01:32.034520+00:00"
$ws.Cells.Item(17, 11).Value = "2025-11-04T07:01:40.359711+00:00"
$ws.Cells.Item(18, 11).Value = "2025-11-04T07:01:43.168626+00:00"
$ws.Cells.Item(19, 11).Value = "2025-11-04T07:01:45.558741+00:00"
$ws.Cells.Item(20, 11).Value = "2025-11-04T07:01:48.459900+00:00"
$ws.Cells.Item(21, 11).Value = "2025-11-04T07:01:48.459933+00:00"
$ws.Cells.Item(22, 11).Value = "2025-11-04T07:01:48.459953+00:00"
$ws.Cells.Item(23, 11).Value = "2025-11-04T07:01:51.229533+00:00"
$ws.Cells.Item(24, 11).Value = "2025-11-04T07:01:51.229565+00:00"
$ws.Cells.Item(25, 11).Value = "2025-11-04T07:01:51.229585+00:00"
$ws.Cells.Item(26, 11).Value = "2025-11-04T07:01:51.229604+00:00"
$ws.Cells.Item(27, 11).Value = "2025-11-04T07:01:53.600381+00:00"
$ws.Cells.Item(28, 11).Value = "2025-11-04T07:01:59.412704+00:00"
$ws.Cells.Item(29, 11).Value = "2025-11-04T07:01:59.412733+00:00"
$ws.Cells.Item(30, 11).Value = "2025-11-04T07:01:59.412751+00:00"
$ws.Cells.Item(31, 11).Value = "2025-11-04T07:01:59.412768+00:00"
$ws.Cells.Item(32, 11).Value = "2025-11-04T07:02:02.190670+00:00"
$ws.Cells.Item(33, 11).Value = "2025-11-04T07:02:02.190700+00:00"
$ws.Cells.Item(34, 11).Value = "2025-11-04T07:02:02.190717+00:00"
$ws.Cells.Item(35, 11).Value = "2025-11-04T07:02:04.547944+00:00"
$ws.Cells.Item(36, 11).Value = "2025-11-04T07:02:04.547974+00:00"
$ws.Cells.Item(37, 11).Value = "2025-11-04T07:02:04.547991+00:00"
$ws.Cells.Item(38, 11).Value = "2025-11-04T07:02:04.548011+00:00"
$ws.Cells.Item(39, 11).Value = "2025-11-04T07:02:04.548028+00:00"
$ws.Cells.Item(40, 11).Value = "2025-11-04T07:02:04.548044+00:00"
$ws.Cells.Item(41, 11).Value = "2025-11-04T07:02:04.548060+00:00"
$ws.Cells.Item(42, 11).Value = "2025-11-04T07:02:04.548074+00:00"
$ws.Cells.Item(43, 11).Value = "2025-11-04T07:02:07.540940+00:00"
$ws.Cells.Item(44, 11).Value = "2025-11-04T07:02:07.540970+00:00"
$ws.Cells.Item(45, 11).Value = "2025-11-04T07:02:12.225726+00:00"
$ws.Cells.Item(46, 11).Value = "2025-11-04T07:02:14.586792+00:00"
$ws.Cells.Item(47, 11).Value = "2025-11-04T07:02:14.586823+00:00"
$ws.Cells.Item(48, 11).Value = "2025-11-04T07:02:14.586842+00:00"
